# Remove the duplicated block of orphan-tag paragraphs.
#
# The document originally contained the orphan-tag list once, then
# repeated the whole list again, then repeated a trailing portion of it
# a third time. This removes the (first) full duplicate copy, which
# spans from the paragraph right after the first "PUMP:URS:4000 " entry
# through the last "PUMP:UNIT:220" entry that immediately precedes the
# point where the (now-trailing) partial repeat resumes with
# "PUMP:HRS:103".

$d = $word.ActiveDocument

$startPara = $d.Paragraphs.Item(18)
$endPara   = $d.Paragraphs.Item(78)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()

Write-Host "Paragraphs remaining:" $d.Paragraphs.Count
